$wb = $excel.ActiveWorkbook

# Update "想去人数" (people interested) counts on the "展览" sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 155
$wsExhibit.Range("F5").Value = 5

# Update the same counts on the "全部类型" sheet (mirrors the same rows)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 155
$wsAll.Range("F5").Value = 5
